$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2194
$ws1.Range("F3").Value = 909
$ws1.Range("F4").Value = 1641
$ws1.Range("F5").Value = 387

# Sheet "全部类型" (sheet4 / rId4) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2194
$ws4.Range("F5").Value = 909
$ws4.Range("F6").Value = 1641
$ws4.Range("F7").Value = 387
